$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - same style as the other headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF) for rows 2-17
$data = @(
    @(2, 4),
    @(8, 8),
    @(4, 6),
    @(6, 7),
    @(7, 7),
    @(7, 8),
    @(9, 9),
    @(6, 7),
    @(5, 7),
    @(6, 7),
    @(10, 10),
    @(5, 7),
    @(3, 4),
    @(4, 6),
    @(4, 6),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
